# Update gh-pages to output generated at 456a3b4
# Applies numeric "want-to-go" count bumps across the four sheets and
# inserts one brand-new event row into the "演出" (shows) sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value  = 563
$ws1.Cells.Item(3, 6).Value  = 39
$ws1.Cells.Item(5, 6).Value  = 6323
$ws1.Cells.Item(6, 6).Value  = 709
$ws1.Cells.Item(9, 6).Value  = 621
$ws1.Cells.Item(10, 6).Value = 305
$ws1.Cells.Item(12, 6).Value = 668
$ws1.Cells.Item(13, 6).Value = 4
$ws1.Cells.Item(14, 6).Value = 1135
$ws1.Cells.Item(15, 6).Value = 75
$ws1.Cells.Item(16, 6).Value = 396
$ws1.Cells.Item(19, 6).Value = 1409
$ws1.Cells.Item(20, 6).Value = 652
$ws1.Cells.Item(21, 6).Value = 368
$ws1.Cells.Item(22, 6).Value = 386
$ws1.Cells.Item(24, 6).Value = 1060
$ws1.Cells.Item(25, 6).Value = 117
$ws1.Cells.Item(26, 6).Value = 2177
$ws1.Cells.Item(27, 6).Value = 237
$ws1.Cells.Item(28, 6).Value = 83
$ws1.Cells.Item(29, 6).Value = 387
$ws1.Cells.Item(31, 6).Value = 3510

# ---------------------------------------------------------------------
# Sheet 2: 演出 (shows)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(4, 6).Value  = 159
$ws2.Cells.Item(5, 6).Value  = 1
$ws2.Cells.Item(6, 7).Value  = 112
$ws2.Cells.Item(7, 6).Value  = 25
$ws2.Cells.Item(8, 6).Value  = 699
$ws2.Cells.Item(24, 6).Value = 181
$ws2.Cells.Item(32, 6).Value = 1581

# Insert a brand-new row 34 ("哈利的魔法世界" candle-light concert) and
# push the old row 34 ("菊次郎的夏天") down to row 35.
$ws2.Rows.Item(34).Insert()

# Copy formatting from the row right above so the new/shifted rows keep
# the same (default) look instead of whatever Excel guesses on insert.
$ws2.Range("A33:I33").Copy()
$ws2.Range("A34:I34").PasteSpecial(-4122)

$ws2.Cells.Item(34, 1).Value = 33

$ws2.Cells.Item(34, 2).NumberFormat = "@"
$ws2.Cells.Item(34, 2).Value = "2024.05.19"
$ws2.Cells.Item(34, 2).Style = "Normal"

$ws2.Cells.Item(34, 3).Value = "上海·《哈利的魔法世界》全系列烛光音乐会"
$ws2.Cells.Item(34, 4).Value = "南京西路1376号 上海商城剧院"

$ws2.Cells.Item(34, 5).NumberFormat = "@"
$ws2.Cells.Item(34, 5).Value = "2024.05.19 19:30-05.19 21:00"
$ws2.Cells.Item(34, 5).Style = "Normal"

$ws2.Cells.Item(34, 6).Value = 1
$ws2.Cells.Item(34, 7).Value = 100
$ws2.Cells.Item(34, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82164"
$ws2.Cells.Item(34, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/OuWJLMSK1709105632932.jpeg"

# The shifted-down row (old row 34, now row 35) keeps its data but its
# index column (A) advances by one, matching the new row number.
$ws2.Cells.Item(35, 1).Value = 34

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 (local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value  = 1814
$ws3.Cells.Item(5, 6).Value  = 1178
$ws3.Cells.Item(7, 6).Value  = 1565
$ws3.Cells.Item(8, 6).Value  = 428
$ws3.Cells.Item(11, 6).Value = 735

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (all types) -- aggregated view of the other sheets,
# updated in lock-step with the same numbers (row layout unaffected by
# the new row 34 inserted above, per the source diff).
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value  = 1814
$ws4.Cells.Item(4, 6).Value  = 1178
$ws4.Cells.Item(5, 6).Value  = 1565
$ws4.Cells.Item(6, 6).Value  = 428
$ws4.Cells.Item(8, 6).Value  = 735
$ws4.Cells.Item(9, 6).Value  = 563
$ws4.Cells.Item(10, 6).Value = 39
$ws4.Cells.Item(12, 6).Value = 6323
$ws4.Cells.Item(13, 7).Value = 112
$ws4.Cells.Item(14, 6).Value = 25
$ws4.Cells.Item(15, 6).Value = 709
$ws4.Cells.Item(17, 6).Value = 699
$ws4.Cells.Item(18, 6).Value = 621
$ws4.Cells.Item(20, 6).Value = 668
$ws4.Cells.Item(25, 6).Value = 1135
$ws4.Cells.Item(26, 6).Value = 396
$ws4.Cells.Item(31, 6).Value = 1409
$ws4.Cells.Item(34, 6).Value = 652
$ws4.Cells.Item(35, 6).Value = 368
$ws4.Cells.Item(36, 6).Value = 386
$ws4.Cells.Item(39, 6).Value = 181
$ws4.Cells.Item(45, 6).Value = 1581
$ws4.Cells.Item(46, 6).Value = 237
$ws4.Cells.Item(47, 6).Value = 83
$ws4.Cells.Item(48, 6).Value = 387
$ws4.Cells.Item(50, 6).Value = 3510
